$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111863073
$ws.Range("B2").Value = 88899
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 3286
$ws.Range("F2").Value = "Flattoppad klubbsvamp"
$ws.Range("G2").Value = "Clavariadelphus truncatus"
$ws.Range("H2").Value = "(Quél.) Donk"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2"
$ws.Range("I2").NumberFormat = "General"
$ws.Range("P2").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q2").Value = 655228
$ws.Range("R2").Value = 6634879
$ws.Range("Z2").Value = "10:50"
$ws.Range("AB2").Value = "10:50"
$ws.Range("Q3").Value = 655138
$ws.Range("R3").Value = 6634821
$ws.Range("A4").Value = 111863001
$ws.Range("B4").Value = 90332
$ws.Range("E4").Value = 4769
$ws.Range("F4").Value = "Svavelriska"
$ws.Range("G4").Value = "Lactarius scrobiculatus"
$ws.Range("H4").Value = "(Scop.:Fr.) Fr."
$ws.Range("P4").Value = "Charlottenberg, Upl"
$ws.Range("Q4").Value = 655218
$ws.Range("R4").Value = 6634940
$ws.Range("Z4").Value = "10:47"
$ws.Range("AB4").Value = "10:47"
$ws.Range("A5").Value = 111863040
$ws.Range("B5").Value = 90687
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 5964
$ws.Range("F5").Value = "Fjällig taggsvamp s.str."
$ws.Range("G5").Value = "Sarcodon imbricatus s.str."
$ws.Range("H5").Value = "(L.:Fr.) P.Karst."
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("Q5").Value = 655235
$ws.Range("R5").Value = 6634878
$ws.Range("Z5").Value = "10:49"
$ws.Range("AB5").Value = "10:49"
$ws.Range("AC5").Value = "Halv häxring, 3 m i diameter"
$ws.Range("Q6").Value = 655135
$ws.Range("R6").Value = 6634793
$ws.Range("A7").Value = 111862959
$ws.Range("B7").Value = 90687
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 5964
$ws.Range("F7").Value = "Fjällig taggsvamp s.str."
$ws.Range("G7").Value = "Sarcodon imbricatus s.str."
$ws.Range("H7").Value = "(L.:Fr.) P.Karst."
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "11"
$ws.Range("I7").NumberFormat = "General"
$ws.Range("P7").Value = "Charlottenberg, Upl"
$ws.Range("Q7").Value = 655218
$ws.Range("R7").Value = 6634940
$ws.Range("Z7").Value = "10:37"
$ws.Range("AB7").Value = "10:37"
$ws.Range("AC7").Value = "Under gran och tall i en svacka"
$ws.Range("A8").Value = 111863269
$ws.Range("B8").Value = 85062
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 249278
$ws.Range("F8").Value = "Barrviolspindling"
$ws.Range("G8").Value = "Cortinarius harcynicus"
$ws.Range("H8").Value = "(Pers.) M.M.Moser"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "4"
$ws.Range("I8").NumberFormat = "General"
$ws.Range("J8").Value = "fruktkroppar"
$ws.Range("Q8").Value = 655135
$ws.Range("R8").Value = 6634800
$ws.Range("Z8").Value = "11:02"
$ws.Range("AB8").Value = "11:02"
$ws.Range("AC8").Value = "4 ex i gräsglänta under gran och tall."
$ws.Range("A9").Value = 111863402
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "1"
$ws.Range("I9").NumberFormat = "General"
$ws.Range("P9").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q9").Value = 655200
$ws.Range("R9").Value = 6634770
$ws.Range("Z9").Value = "11:02"
$ws.Range("AB9").Value = "11:02"
$ws.Range("AC9").Value = ""
$ws.Range("Q10").Value = 655234
$ws.Range("R10").Value = 6634889
